# Add a new "metadata_schema_id" column (S) to the "Export as TSV" sheet,
# mirroring the formatting of the last existing header column (R), and
# attach the explanatory header comment - matching the CEDAR docs update
# described in the commit ("Update with source murine CEDAR docs - #52").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastHeader = $ws.Range("R1")
$newHeader  = $ws.Range("S1")

# Copy the header cell's formatting (bold, centered, wrapped text, etc.)
# onto the new column before filling in its value/comment.
$lastHeader.Copy()
$newHeader.PasteSpecial(-4122) # xlPasteFormats

$newHeader.Value = "metadata_schema_id"
$newHeader.AddComment("The string that serves as the definitive identifier for the metadata schema version.")

# Keep the visible selection where the author left it after adding the column.
$ws.Range("S3").Select() | Out-Null

Write-Output "Added metadata_schema_id column with header comment"
